$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two remaining channel-requirement columns (header row 1):
#   F1: requirements-print     -> requirements-GS1_GDSN
#   G1: requirements-ecommerce -> requirements-PRODUCT_CATALOG
$ws.Range("F1").Value = "requirements-GS1_GDSN"
$ws.Range("G1").Value = "requirements-PRODUCT_CATALOG"

# The third channel column (H: requirements-mobile / sku placeholders) was
# dropped entirely, so remove column H and shift everything left.
$ws.Columns.Item(8).Delete()

# Mirror the saved view state (selection moved to A2).
$ws.Range("A2").Select() | Out-Null
